$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new review rows to the log. Cell-write order below is chosen
# so that new shared-string entries are interned in the same order as the
# target workbook (first occurrence order matters for xl/sharedStrings.xml).
$ws.Range("A8").Value = "RV_00_70"
$ws.Range("B8").Value = "LM35 CDD"
$ws.Range("D8").Value = "MED"
$ws.Range("E8").Value = "Kareem"

$ws.Range("A9").Value = "RV_00_80"
$ws.Range("B9").Value = "Software"

$ws.Range("E10").Value = "Alhassan"

$ws.Range("C9").Value = "Need to be organized "
$ws.Range("C10").Value = "missing some requirements"

$ws.Range("A10").Value = "RV_00_90"

# Fill in the remaining cells, which reuse already-existing shared strings.
$ws.Range("C8").Value = "ADDING Some requirements "
$ws.Range("D9").Value = "Med"
$ws.Range("E9").Value = "Kareem"
$ws.Range("B10").Value = "Software"
$ws.Range("D10").Value = "Med"

# Widen column C so the longer "Changes" text fits (closest reachable value
# to the source workbook's best-fit width given this engine's column-width
# quantization).
$ws.Columns.Item(3).ColumnWidth = 46.5

# Match the active-cell selection left behind in the source workbook.
$ws.Range("K11").Select() | Out-Null
